$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 3
$ws.Range("H2").Value = 3.25
$ws.Range("I2").Value = 2.3
$ws.Range("J2").Value = 3.75
$ws.Range("L2").Value = 3.1
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 10
$ws.Range("Q2").Value = 2.1
$ws.Range("R2").Value = 1.7
$ws.Range("S2").Value = 3.75
$ws.Range("T2").Value = 1.25
$ws.Range("AA2").Value = 11
$ws.Range("AC2").Value = 26
$ws.Range("AK2").Value = 11
